$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# ------------------------------------------------------------------
# 1. Capture formatting from rows that are about to be overwritten
#    BEFORE we overwrite them (order matters).
# ------------------------------------------------------------------

# old row 100 (summary row) B:G  ->  new row 105 B:G
$ws.Range("B100:G100").Copy()
$ws.Range("B105:G105").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# old row 99 (blank row, F:K styled) F:K  ->  new row 103 F:K
$ws.Range("F99:K99").Copy()
$ws.Range("F103:K103").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# old row 99 F:G  ->  new row 104 F:G
$ws.Range("F99:G99").Copy()
$ws.Range("F104:G104").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# old row 99 J:K  ->  new row 102 J:K
$ws.Range("J99:K99").Copy()
$ws.Range("J102:K102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Apply the "data row" formatting (taken from row 98) to the three
#    new data rows 99, 100, 101.
# ------------------------------------------------------------------
$ws.Range("A98:K98").Copy()
$ws.Range("A99:K99").PasteSpecial(-4122)
$ws.Range("A100:K100").PasteSpecial(-4122)
$ws.Range("A101:K101").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3. Fill in the new data.
# ------------------------------------------------------------------

# Row 99
$ws.Range("A99").Value = 21
$ws.Range("B99").Value = "Interface Design"
$ws.Range("C99").Value = "Pseudo Code"
$ws.Range("D99").Value = "[TASK]"
$ws.Range("E99").Value = "Erklaerung der Struktur"
$ws.Range("F99").Value = 44369
$ws.Range("G99").Value = 44359
$ws.Range("I99").Formula = "=ROUNDUP(((SUM(K99-J99)*24*60/60)/0.25),0)*0.25"
$ws.Range("J99").Value = 0.72916666666666663
$ws.Range("K99").Value = 0.75

# Row 100
$ws.Range("A100").Value = 22
$ws.Range("B100").Value = "Interface Design"
$ws.Range("C100").Value = "Pseudo Code"
$ws.Range("D100").Value = "[FEATURE]"
$ws.Range("E100").Value = "Pseudocode fuer Aehnliche Rezepte"
$ws.Range("F100").Value = 44372
$ws.Range("G100").Value = 44359
$ws.Range("I100").Formula = "=ROUNDUP(((SUM(K100-J100)*24*60/60)/0.25),0)*0.25"
$ws.Range("J100").Value = 0.58333333333333337
$ws.Range("K100").Value = 0.66666666666666663

# Row 101
$ws.Range("A101").Value = 23
$ws.Range("B101").Value = "Interface Design"
$ws.Range("C101").Value = "Pseudo Code"
$ws.Range("D101").Value = "[FEATURE]"
$ws.Range("E101").Value = "Pseudocode fuer Empfohlene Rezepte"
$ws.Range("F101").Value = 44372
$ws.Range("G101").Value = 44359
$ws.Range("I101").Formula = "=ROUNDUP(((SUM(K101-J101)*24*60/60)/0.25),0)*0.25"
$ws.Range("J101").Value = 0.66666666666666663
$ws.Range("K101").Value = 0.75

# ------------------------------------------------------------------
# 4. Re-write the summary row at its new location (row 105). The
#    formulas stay the same (they already sum whole columns), so we
#    just need to re-enter them in the new row.
# ------------------------------------------------------------------
$ws.Range("B105").Value = "Stunden insgesamt"
$ws.Range("C105").Formula = "=SUM(I:I)+SUM(H:H)"
$ws.Range("D105").Value = "Stunden Seminar"
$ws.Range("E105").Formula = "=SUM(H:H)"
$ws.Range("F105").Value = "Stunden Projekt"
$ws.Range("G105").Formula = "=SUM(I:I)"

# ------------------------------------------------------------------
# 5. Data validation: the "D" prefix dropdown needs to keep covering
#    the data rows (now through D101) and the new blank rows
#    (D103:D104).
# ------------------------------------------------------------------
$ws.Range("D41:D99").Validation.Delete()

$rNew1 = $ws.Range("D41:D101")
$rNew1.Validation.Add(3, 1, 1, "=`$N`$3:`$N`$6")
$rNew1.Validation.ErrorTitle = "Prefix nicht unterstützt"
$rNew1.Validation.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden`n"
$rNew1.Validation.InputTitle = "Prefix"
$rNew1.Validation.InputMessage = "Wählen Sie einen Prefix aus"

$rNew2 = $ws.Range("D103:D104")
$rNew2.Validation.Add(3, 1, 1, "=`$N`$3:`$N`$6")
$rNew2.Validation.ErrorTitle = "Prefix nicht unterstützt"
$rNew2.Validation.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden`n"
$rNew2.Validation.InputTitle = "Prefix"
$rNew2.Validation.InputMessage = "Wählen Sie einen Prefix aus"

# ------------------------------------------------------------------
# 6. Selection / view state (best effort).
# ------------------------------------------------------------------
$ws.Range("J102").Select()

Write-Host "Done"
